$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move A4's value up into A3, then clear A4. This reproduces the diff:
# A3 goes from shared string "RO.ACT.001.MAJ.01" to "AD.SEC.002.FON.01"
# (the value that used to live in A4), A4 becomes empty, and the now
# unreferenced shared string "RO.ACT.001.MAJ.01" drops out of the
# sharedStrings table.
$v = $ws.Range("A4").Value()
$ws.Range("A3").Value = $v
$ws.Range("A4").ClearContents()

# Restore the selection recorded in the saved view state (A8:A9).
$ws.Range("A8:A9").Select()
